# On Pilgrimage - October 1957 (DDLW #729)
# Convert the old "Heading1 + bold byline" title block into a pandoc-style
# title block: a Title-styled paragraph followed by an Authors-styled
# paragraph, with the wrapping bookmark removed and each word/space of the
# title broken into its own run (matching the source markdown -> docx
# conversion's run-per-inline-token layout).

$d = $word.ActiveDocument

# --- 1. Remove the old "On Pilgrimage - October 1957" Heading1 paragraph --
# (it also happens to be wrapped in a bookmark). Deleting the paragraph's
# range leaves the now-empty bookmarkStart/bookmarkEnd pair touching at
# position 0; two subsequent zero-length deletes at that position peel
# off each marker in turn without touching any real text.
$d.Paragraphs(1).Range.Delete()
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- 2. Insert a new empty paragraph in front of the byline paragraph ----
# (currently paragraph 1 is "By Dorothy Day") to hold the new Title.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphBefore() | Out-Null
$titlePara = $d.Paragraphs(1)

$titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:pPr><w:pStyle w:val='Title'/></w:pPr>" +
  "<w:r><w:t xml:space='preserve'>On</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>Pilgrimage</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>-</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>October</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>1957</w:t></w:r>" +
  "</w:p>"
$titlePara.Range.InsertXML($titleXml)

# --- 3. Replace "By Dorothy Day" (bold, unstyled) with an Authors-styled -
# paragraph reading just "Dorothy Day".
$authorsPara = $d.Paragraphs(2)
$authorsXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:pPr><w:pStyle w:val='Authors'/></w:pPr>" +
  "<w:r><w:t xml:space='preserve'>Dorothy</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>Day</w:t></w:r>" +
  "</w:p>"
$authorsPara.Range.InsertXML($authorsXml)
